$wb = $excel.ActiveWorkbook

# Rename sheets p3..p9 to p03..p09 (zero-padded to match p01, p02, p10-p13 naming)
$renames = @{
    "p3" = "p03"
    "p4" = "p04"
    "p5" = "p05"
    "p6" = "p06"
    "p7" = "p07"
    "p8" = "p08"
    "p9" = "p09"
}
foreach ($old in $renames.Keys) {
    $ws = $wb.Worksheets.Item($old)
    $ws.Name = $renames[$old]
}

# Previously the active/selected tab was p03 (formerly "p3", the 4th sheet) with
# cell B10 selected. Move the active tab to p13 (the last sheet) and select B15 there,
# leaving p03's own remembered selection on B10.
$p03 = $wb.Worksheets.Item("p03")
[void]$p03.Range("B10").Select()

$p13 = $wb.Worksheets.Item("p13")
[void]$p13.Select()
[void]$p13.Range("B15").Select()
